$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '64.121.08'
Set-TextCell $ws 'E2' '  -4.76%  '

Set-TextCell $ws 'D3' '3.294.24'
Set-TextCell $ws 'E3' '  -6.08%  '

Set-TextCell $ws 'D4' '0.999'
Set-TextCell $ws 'E4' '  -0.25%  '

Set-TextCell $ws 'D5' '178.28'
Set-TextCell $ws 'E5' '  -11.51%  '

Set-TextCell $ws 'D6' '525.12'
Set-TextCell $ws 'E6' '  -5.19%  '

Set-TextCell $ws 'D7' '0.605'
Set-TextCell $ws 'E7' '  -0.53%  '

Set-TextCell $ws 'D8' '3.291.17'
Set-TextCell $ws 'E8' '  -5.91%  '

Set-TextCell $ws 'E9' '  -0.05%  '

Set-TextCell $ws 'D10' '0.609'
Set-TextCell $ws 'E10' '  -7.09%  '

Set-TextCell $ws 'D11' '57.66'
Set-TextCell $ws 'E11' '  -8.02%  '

Set-TextCell $ws 'D12' '0.133'
Set-TextCell $ws 'E12' '  -7.04%  '

Set-TextCell $ws 'D13' '0.0000259'
Set-TextCell $ws 'E13' '  -4.52%  '

Set-TextCell $ws 'D14' '9.12'
Set-TextCell $ws 'E14' '  -7.43%  '

Set-TextCell $ws 'D15' '3.803.91'
Set-TextCell $ws 'E15' '  -6.49%  '

Set-TextCell $ws 'B16' 'TRON'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D16' '0.117'
Set-TextCell $ws 'E16' '  -5.28%  '

Set-TextCell $ws 'B17' 'WrappedEther'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 'D17' '3.281.20'
Set-TextCell $ws 'E17' '  -6.44%  '

Set-TextCell $ws 'D18' '63.915.42'
Set-TextCell $ws 'E18' '  -4.70%  '

Set-TextCell $ws 'D19' '17.47'
Set-TextCell $ws 'E19' '  -5.35%  '

Set-TextCell $ws 'D20' '11.12'
Set-TextCell $ws 'E20' '  -6.05%  '

Set-TextCell $ws 'D21' '0.958'
Set-TextCell $ws 'E21' '  -6.80%  '

Set-TextCell $ws 'D22' '374.78'
Set-TextCell $ws 'E22' '  -4.43%  '

Set-TextCell $ws 'D23' '3.77'
Set-TextCell $ws 'E23' '  -5.89%  '

Set-TextCell $ws 'D24' '80.73'
Set-TextCell $ws 'E24' '  -2.77%  '

Set-TextCell $ws 'D25' '11.09'
Set-TextCell $ws 'E25' '  -11.69%  '

Set-TextCell $ws 'D26' '3.90'
Set-TextCell $ws 'E26' '  -0.75%  '

Set-TextCell $ws 'B27' 'LEO'
Set-TextCell $ws 'C27' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws 'D27' '6.08'
Set-TextCell $ws 'E27' '  -1.55%  '

Set-TextCell $ws 'B28' 'ImmutableX'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D28' '2.68'
Set-TextCell $ws 'E28' '  -4.94%  '

Set-TextCell $ws 'B29' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D29' '11.41'
Set-TextCell $ws 'E29' '  -7.08%  '

Set-TextCell $ws 'B30' 'Filecoin'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D30' '8.37'
Set-TextCell $ws 'E30' '  -5.53%  '

Set-TextCell $ws 'B31' 'EthereumClassic'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D31' '28.90'
Set-TextCell $ws 'E31' '  -6.88%  '

Set-TextCell $ws 'B32' 'Bittensor'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D32' '638.35'
Set-TextCell $ws 'E32' '  -6.97%  '

Set-TextCell $ws 'B33' 'NEARProtocol'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D33' '6.66'
Set-TextCell $ws 'E33' '  -5.48%  '

Set-TextCell $ws 'B34' 'Cosmos'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D34' '11.27'
Set-TextCell $ws 'E34' '  -3.94%  '

Set-TextCell $ws 'D35' '59.35'
Set-TextCell $ws 'E35' '  -6.93%  '

Set-TextCell $ws 'B36' 'Hedera'
Set-TextCell $ws 'C36' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D36' '0.106'
Set-TextCell $ws 'E36' '  -5.49%  '

Set-TextCell $ws 'B37' 'Dai'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D37' '1.00'
Set-TextCell $ws 'E37' '  +0.00%  '

Set-TextCell $ws 'B38' 'TheGraph'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws 'D38' '0.390'
Set-TextCell $ws 'E38' '  -2.21%  '

Set-TextCell $ws 'B39' 'InjectiveProtocol'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D39' '36.70'
Set-TextCell $ws 'E39' '  -5.33%  '

Set-TextCell $ws 'B40' 'FirstDigitalUSD'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws 'D40' '0.995'
Set-TextCell $ws 'E40' '  -0.27%  '

Set-TextCell $ws 'B41' 'PEPE'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws 'D41' '0.0₃0704'
Set-TextCell $ws 'E41' '  +4.15%  '

Set-TextCell $ws 'B42' 'Maker'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D42' '2.941.20'
Set-TextCell $ws 'E42' '  -4.32%  '

Set-TextCell $ws 'B43' 'Kaspa'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D43' '0.124'
Set-TextCell $ws 'E43' '  -4.99%  '

Set-TextCell $ws 'B44' 'Fetch.AI'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D44' '2.46'
Set-TextCell $ws 'E44' '  -6.14%  '

Set-TextCell $ws 'B45' 'ThetaToken'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell $ws 'D45' '2.70'
Set-TextCell $ws 'E45' '  -10.00%  '

Set-TextCell $ws 'D46' '0.0399'
Set-TextCell $ws 'E46' '  -0.62%  '

Set-TextCell $ws 'B47' 'WEMIXToken'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws 'D47' '2.64'
Set-TextCell $ws 'E47' '  -4.83%  '

Set-TextCell $ws 'B48' 'ApeXProtocol'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell $ws 'D48' '3.01'
Set-TextCell $ws 'E48' '  +4.44%  '

Set-TextCell $ws 'B49' 'Stacks'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D49' '2.78'
Set-TextCell $ws 'E49' '  +5.79%  '

Set-TextCell $ws 'B50' 'Stellar'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D50' '0.126'
Set-TextCell $ws 'E50' '  -1.17%  '

Set-TextCell $ws 'B51' 'Monero'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D51' '135.87'
Set-TextCell $ws 'E51' '  -1.65%  '
